$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D4").Value = 23252
$ws.Range("H10").Value = 491483
$ws.Range("I10").Value = 28893
$ws.Range("H11").Value = 401376
$ws.Range("I11").Value = 53820
$ws.Range("H12").Value = 84617
$ws.Range("I12").Value = 5879
$ws.Range("H13").Value = 311040
$ws.Range("I13").Value = 41325
$ws.Range("H14").Value = 288383
$ws.Range("I14").Value = 45649
$ws.Range("H15").Value = 182978
$ws.Range("I15").Value = 49256
$ws.Range("H16").Value = 182453
$ws.Range("I16").Value = 50149
$ws.Range("H17").Value = 383594
$ws.Range("I17").Value = 52058
$ws.Range("H18").Value = 223230
$ws.Range("I18").Value = 58611
$ws.Range("H19").Value = 228744
$ws.Range("I19").Value = 59500
$ws.Range("H20").Value = 322422
$ws.Range("I20").Value = 64703
$ws.Range("H21").Value = 338294
$ws.Range("I21").Value = 87634
$ws.Range("H22").Value = 535311
$ws.Range("I22").Value = 89011
$ws.Range("H23").Value = 479548
$ws.Range("I23").Value = 114460
$ws.Range("H24").Value = 634760
$ws.Range("I24").Value = 115618
$ws.Range("H25").Value = 48864
$ws.Range("I25").Value = 14541
$ws.Range("H26").Value = 677698
$ws.Range("I26").Value = 163394
$ws.Range("H27").Value = 431847
$ws.Range("I27").Value = 175531
$ws.Range("H28").Value = 701813
$ws.Range("I28").Value = 191605
$ws.Range("H29").Value = 785539
$ws.Range("I29").Value = 197300
$ws.Range("H30").Value = 634824
$ws.Range("I30").Value = 215229
$ws.Range("H31").Value = 725075
$ws.Range("I31").Value = 233456
$ws.Range("H32").Value = 884752
$ws.Range("I32").Value = 295780
$ws.Range("H33").Value = 22825
$ws.Range("I33").Value = 2200
$ws.Range("H36").Value = 27205
$ws.Range("I36").Value = 6033
$ws.Range("H37").Value = 29822
$ws.Range("H38").Value = 13037
$ws.Range("I38").Value = 2666
$ws.Range("H39").Value = 19493
$ws.Range("I39").Value = 2681
$ws.Range("H41").Value = 10867
$ws.Range("I41").Value = 3806
$ws.Range("H42").Value = 36519
$ws.Range("I42").Value = 5811
$ws.Range("H43").Value = 32105
$ws.Range("I43").Value = 7068
$ws.Range("H44").Value = 27259
$ws.Range("I44").Value = 7105
$ws.Range("H45").Value = 35637
$ws.Range("I45").Value = 8782
$ws.Range("H46").Value = 17881
$ws.Range("I46").Value = 1766
$ws.Range("H47").Value = 51264
$ws.Range("I47").Value = 9953
$ws.Range("H49").Value = 42052
$ws.Range("I49").Value = 4644
$ws.Range("H50").Value = 43266
$ws.Range("I50").Value = 10504
$ws.Range("H51").Value = 56109
$ws.Range("I51").Value = 10540
$ws.Range("H52").Value = 62370
$ws.Range("I52").Value = 10777
$ws.Range("H54").Value = 40827
$ws.Range("I54").Value = 6527
$ws.Range("H55").Value = 30892
$ws.Range("I55").Value = 2634
$ws.Range("H56").Value = 37060
$ws.Range("I56").Value = 8327
$ws.Range("H59").Value = 26986
$ws.Range("H60").Value = 32526
$ws.Range("I60").Value = 4734
$ws.Range("H62").Value = 62831
$ws.Range("I62").Value = 16448
$ws.Range("H65").Value = 34364
$ws.Range("H67").Value = 28299
$ws.Range("I67").Value = 6628
$ws.Range("H68").Value = 35078
$ws.Range("I68").Value = 5885
$ws.Range("H70").Value = 183422
$ws.Range("I70").Value = 15346
$ws.Range("H71").Value = 85496
$ws.Range("I71").Value = 7494
$ws.Range("H72").Value = 58953
$ws.Range("I72").Value = 15417
$ws.Range("H73").Value = 144265
$ws.Range("I73").Value = 18325
$ws.Range("H74").Value = 108982
$ws.Range("I74").Value = 30880
$ws.Range("H75").Value = 451114
$ws.Range("I75").Value = 66022
$ws.Range("H76").Value = 121381
$ws.Range("I76").Value = 53079
$ws.Range("H77").Value = 322677
$ws.Range("I77").Value = 32362
$ws.Range("H78").Value = 462151
$ws.Range("I78").Value = 57540
$ws.Range("H79").Value = 276785
$ws.Range("I79").Value = 99199
$ws.Range("H81").Value = 380343
$ws.Range("I81").Value = 38674
$ws.Range("H82").Value = 235945
$ws.Range("I82").Value = 54347
$ws.Range("H83").Value = 312313
$ws.Range("I83").Value = 110456
$ws.Range("H84").Value = 712270
$ws.Range("I84").Value = 191921
$ws.Range("H85").Value = 170184
$ws.Range("I85").Value = 39831
$ws.Range("H86").Value = 428403
$ws.Range("I86").Value = 59051
$ws.Range("H87").Value = 289159
$ws.Range("I87").Value = 62059
$ws.Range("H88").Value = 325040
$ws.Range("I88").Value = 58382
$ws.Range("H89").Value = 184202
$ws.Range("I89").Value = 67663
$ws.Range("H90").Value = 313731
$ws.Range("I90").Value = 115692
$ws.Range("H91").Value = 503295
$ws.Range("I91").Value = 141103
$ws.Range("H92").Value = 278903
$ws.Range("I92").Value = 124849
$ws.Range("H93").Value = 395520
$ws.Range("I93").Value = 102394
$ws.Range("H94").Value = 557226
$ws.Range("I94").Value = 230436
$ws.Range("H95").Value = 6448310
$ws.Range("I95").Value = 5328194

Write-Output "Applied 144 cell updates"
